# Excel COM-interop script to apply the EPEX spot prices workbook update.
#
# Summary of changes:
# 1. "Prix Spot" sheet: insert a new column before column DU (the 125th
#    column). The new column gets header "20-nov" (row 1) and "-" for
#    rows 2-25 (same pattern as the other "no data" placeholder columns).
#    This shifts every column from the old DU onward one position to the
#    right (DU->DV, DV->DW, ..., EY->EZ), which Excel does automatically
#    when inserting a column.
# 2. "Gaz" sheet: append a new row 154 with Date=2025-11-18 and
#    Last Price=30.615.
# 3. "CO2" sheet: append a new row 154 with Date=2025-11-18 and
#    Last Price=80.93000000000001.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Prix Spot: insert new column DU ("20-nov") shifting existing data
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Columns("DU:DU").Insert()

$wsPrix.Range("DU1").Value = "20-nov"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 125).Value = "-"
}

# ---------------------------------------------------------------------
# 2. Gaz: append row 154
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A154").Value = "'2025-11-18"
$wsGaz.Range("A154").Style = "Normal"
$wsGaz.Range("B154").Value = 30.615

# ---------------------------------------------------------------------
# 3. CO2: append row 154
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A154").Value = "'2025-11-18"
$wsCO2.Range("A154").Style = "Normal"
$wsCO2.Range("B154").Value = 80.93000000000001
